# Apply the crypto-price refresh described by the commit:
# updates the Price (D) and Volume(1h) (E) columns for the coin rows,
# including a row-49/row-50 swap (Algorand <-> EnergySwap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value must be
# forced to Text (Excel auto-converts single-dot numeric-looking strings
# like "219.70" or "7.63" to numbers, which would lose the original
# text formatting / trailing zeros that the source data uses).
$updates = @(
    @{ Cell = "D2"; Value = "27.133.37"; ForceText = $true }
    @{ Cell = "E2"; Value = "  +1.05%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "1.649.74"; ForceText = $true }
    @{ Cell = "E3"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  -0.84%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "219.70"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.17%  "; ForceText = $false }
    @{ Cell = "E6"; Value = "  -0.53%  "; ForceText = $false }
    @{ Cell = "E7"; Value = "  -0.76%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.63%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "19.70"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +2.19%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  +0.28%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "1.880.52"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "1.642.78"; ForceText = $true }
    @{ Cell = "E13"; Value = "  -0.16%  "; ForceText = $false }
    @{ Cell = "E14"; Value = "  +0.61%  "; ForceText = $false }
    @{ Cell = "E15"; Value = "  +0.62%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "66.27"; ForceText = $true }
    @{ Cell = "E16"; Value = "  +2.09%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "27.105.35"; ForceText = $true }
    @{ Cell = "E17"; Value = "  +0.97%  "; ForceText = $false }
    @{ Cell = "E18"; Value = "  -0.39%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "221.81"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +3.05%  "; ForceText = $false }
    @{ Cell = "E20"; Value = "  -0.87%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "6.80"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +8.35%  "; ForceText = $false }
    @{ Cell = "E22"; Value = "  +0.64%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "2.42"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -2.77%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "147.42"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "E26"; Value = "  -0.76%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "7.40"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +2.62%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  +0.07%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "15.91"; ForceText = $true }
    @{ Cell = "E29"; Value = "  +1.32%  "; ForceText = $false }
    @{ Cell = "E30"; Value = "  +0.35%  "; ForceText = $false }
    @{ Cell = "E31"; Value = "  +1.01%  "; ForceText = $false }
    @{ Cell = "E32"; Value = "  +0.62%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  +0.07%  "; ForceText = $false }
    @{ Cell = "E34"; Value = "  +2.25%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "1.269.13"; ForceText = $true }
    @{ Cell = "E35"; Value = "  -2.11%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "2.43"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -0.64%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  -1.61%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.538"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.03%  "; ForceText = $false }
    @{ Cell = "E39"; Value = "  +0.34%  "; ForceText = $false }
    @{ Cell = "E40"; Value = "  -0.77%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  -0.01%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  +0.69%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "1.790.26"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +0.29%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "62.01"; ForceText = $true }
    @{ Cell = "E44"; Value = "  -0.11%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "92.63"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +0.71%  "; ForceText = $false }
    @{ Cell = "E46"; Value = "  -7.66%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  -0.09%  "; ForceText = $false }
    @{ Cell = "E48"; Value = "  -1.07%  "; ForceText = $false }
    @{ Cell = "B49"; Value = "EnergySwap"; ForceText = $false }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; ForceText = $false }
    @{ Cell = "D49"; Value = "7.63"; ForceText = $true }
    @{ Cell = "E49"; Value = "  -0.56%  "; ForceText = $false }
    @{ Cell = "B50"; Value = "Algorand"; ForceText = $false }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; ForceText = $false }
    @{ Cell = "D50"; Value = "0.0975"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "E51"; Value = "  -0.40%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Stash original format, write as Text so Excel keeps the literal
        # string (no numeric coercion / exponent / rounding), then restore
        # the cell's formatting so no stray style is left behind.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
